$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 6653.2
$ws.Range("I40").Value = 4579.8
$ws.Range("K40").Value = 4579.8
$ws.Range("M40").Value = -4404.8

$ws.Range("H51").Value = 6300.6
$ws.Range("I51").Value = 3500
$ws.Range("J51").Value = 7000.75
$ws.Range("K51").Value = 3500
$ws.Range("L51").Value = 7000.75
$ws.Range("M51").Value = -3016
$ws.Range("N51").Value = -7968.75

$ws.Range("H98").Value = 3994.1
$ws.Range("I98").Value = 984.6667
$ws.Range("K98").Value = 984.6667
$ws.Range("M98").Value = 513.3333

$ws.Range("H106").Value = 2874.3157
$ws.Range("I106").Value = 3183.6365
$ws.Range("K106").Value = 3183.6365
$ws.Range("M106").Value = -2552.6365

$ws.Range("H122").Value = 3994.1
$ws.Range("I122").Value = 984.6667
$ws.Range("K122").Value = 2954.0001
$ws.Range("M122").Value = -504.0001000000002

$ws.Range("H138").Value = 5204.868
$ws.Range("I138").Value = 3646.5625
$ws.Range("K138").Value = 10939.6875
$ws.Range("M138").Value = -5799.6875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3476.8276
$ws.Range("I32").Value = 3038.2593
$ws.Range("J32").Value = 9397.5
$ws.Range("K32").Value = 3038.2593
$ws.Range("L32").Value = 9397.5
$ws.Range("M32").Value = -2751.2593
$ws.Range("N32").Value = -9971.5

$ws.Range("H45").Value = 5516
$ws.Range("I45").Value = 3536.7
$ws.Range("J45").Value = 7165.4165
$ws.Range("K45").Value = 3536.7
$ws.Range("L45").Value = 7165.4165
$ws.Range("M45").Value = -3159.7
$ws.Range("N45").Value = -7919.4165

$ws.Range("H95").Value = 19000
$ws.Range("J95").Value = 19000
$ws.Range("L95").Value = 19000
$ws.Range("N95").Value = -24492

$ws.Range("H122").Value = 4432.0435
$ws.Range("I122").Value = 2991.5386
$ws.Range("K122").Value = 8974.6158
$ws.Range("M122").Value = -6524.6158

$ws.Range("H132").Value = 6507.8125
$ws.Range("I132").Value = 2312.75
$ws.Range("K132").Value = 6938.25
$ws.Range("M132").Value = -4408.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3666.2778
$ws.Range("I99").Value = 3469.7334
$ws.Range("K99").Value = 3469.7334
$ws.Range("M99").Value = -1971.7334

$ws.Range("H105").Value = 4377.3213
$ws.Range("I105").Value = 2503.2
$ws.Range("K105").Value = 2503.2
$ws.Range("M105").Value = -756.1999999999998

$ws.Range("H107").Value = 1773.9259
$ws.Range("I107").Value = 1072.9615
$ws.Range("K107").Value = 1072.9615
$ws.Range("M107").Value = 847.0385000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 2000
$ws.Range("I36").Value = 2000
$ws.Range("K36").Value = 2000
$ws.Range("M36").Value = -1612

$ws.Range("H40").Value = 2000
$ws.Range("I40").Value = 2000
$ws.Range("K40").Value = 2000
$ws.Range("M40").Value = -1840

$ws.Range("H58").Value = 316129.7
$ws.Range("I58").Value = 667686.25
$ws.Range("K58").Value = 667686.25
$ws.Range("M58").Value = -667483.25

$ws.Range("H74").Value = 65157
$ws.Range("J74").Value = 65157
$ws.Range("L74").Value = 65157
$ws.Range("N74").Value = -66905

$ws.Range("H77").Value = 65157
$ws.Range("J77").Value = 65157
$ws.Range("L77").Value = 195471
$ws.Range("N77").Value = -204207

$ws.Range("H86").Value = 55561.625
$ws.Range("I86").Value = 82898.60000000001
$ws.Range("K86").Value = 82898.60000000001
$ws.Range("M86").Value = -81775.60000000001

$ws.Range("H89").Value = 55561.625
$ws.Range("I89").Value = 82898.60000000001
$ws.Range("K89").Value = 414493
$ws.Range("M89").Value = -408877

$ws.Range("H94").Value = 817.9
$ws.Range("J94").Value = 964.4
$ws.Range("L94").Value = 964.4
$ws.Range("N94").Value = -1866.4

$ws.Range("H136").Value = 316129.7
$ws.Range("I136").Value = 667686.25
$ws.Range("K136").Value = 2003058.75
$ws.Range("M136").Value = -2000508.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 207.14285
$ws.Range("I12").Value = 1
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 170

$ws.Range("H70").Value = 950
$ws.Range("I70").Value = 950
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 2850
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -2535
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 950
$ws.Range("I73").Value = 950
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 2850
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -1758
$ws.Range("N73").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 361257.75
$ws.Range("I132").Value = 403708.75
$ws.Range("K132").Value = 1211126.25
$ws.Range("M132").Value = -1208596.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1936.1818
$ws.Range("I22").Value = 1657.4286
$ws.Range("J22").Value = 2424
$ws.Range("K22").Value = 1657.4286
$ws.Range("L22").Value = 2424
$ws.Range("M22").Value = -1362.4286
$ws.Range("N22").Value = -3014

$ws.Range("H27").Value = 1936.1818
$ws.Range("I27").Value = 1657.4286
$ws.Range("J27").Value = 2424
$ws.Range("K27").Value = 1657.4286
$ws.Range("L27").Value = 2424
$ws.Range("M27").Value = -1550.4286
$ws.Range("N27").Value = -2638

$ws.Range("H46").Value = 3265.2666
$ws.Range("J46").Value = 4223.222
$ws.Range("L46").Value = 4223.222
$ws.Range("N46").Value = -4599.222

$ws.Range("H55").Value = 908.8182
$ws.Range("I55").Value = 553.1667
$ws.Range("J55").Value = 1335.6
$ws.Range("K55").Value = 553.1667
$ws.Range("L55").Value = 1335.6
$ws.Range("M55").Value = -380.1667
$ws.Range("N55").Value = -1681.6

$ws.Range("H130").Value = 80429
$ws.Range("J130").Value = 80429
$ws.Range("L130").Value = 80429
$ws.Range("N130").Value = -90469

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 7999.5
$ws.Range("J69").Value = 7999.5
$ws.Range("L69").Value = 7999.5
$ws.Range("N69").Value = -9497.5

$ws.Range("H72").Value = 7999.5
$ws.Range("J72").Value = 7999.5
$ws.Range("L72").Value = 23998.5
$ws.Range("N72").Value = -31486.5

Write-Output "Applied 34 row updates across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR"
